$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "Asset Data" before "Differentiators" (so it
#    becomes the 2nd sheet, right after README).
# ---------------------------------------------------------------------------
$assetData = $wb.Worksheets.Add($wb.Worksheets.Item("Differentiators"))
$assetData.Name = "Asset Data"

# ---------------------------------------------------------------------------
# 2. Populate the "Asset Data" sheet.
# ---------------------------------------------------------------------------
$assetData.Range("A1:F2").Font.Bold = $true

$assetData.Range("A1").Value = "data_model"
$assetData.Range("B1").Value = "mapping"

$assetData.Range("A2").Value = "source"
$assetData.Range("B2").Value = "condition"
$assetData.Range("F2").Value = "differentiator"

$assetData.Range("A3").Value = "asset_model"
$assetData.Range("B3").Value = "t"
$assetData.Range("C3").Value = "condition.perfect"
$assetData.Range("D3").Value = "external_diameter"
$assetData.Range("E3").Value = "wall_thickness"
$assetData.Range("F3").Value = "material"
$assetData.Range("G3").Value = "treatment"

$assetData.Range("A4").Value = "asset_data"
$assetData.Range("B4").Value = "age"
$assetData.Range("D4").Value = "agd"

# Column widths (A:E) to match the source workbook formatting.
$assetData.Range("A1:E1").ColumnWidth = 16.67

# ---------------------------------------------------------------------------
# 3. Update per-sheet selections (cursor position) that changed in the diff.
# ---------------------------------------------------------------------------
$assetData.Activate()
$assetData.Range("C15").Select()

$differentiators = $wb.Worksheets.Item("Differentiators")
$differentiators.Activate()
$differentiators.Range("D22").Select()

$modelParameters = $wb.Worksheets.Item("Model Parameters")
$modelParameters.Activate()
$modelParameters.Range("C8").Select()

$simple = $wb.Worksheets.Item("Simple")
$simple.Activate()
$simple.Range("G17").Select()
